# Auto-generated: updates cryptos list Price (D) / Volume(1h) (E) columns
# per the commit "Updated cryptos list on Tue May  2 10:46:27 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.060.23"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "1.831.95"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"

$ws.Range("E5").Value = "  -3.62%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4613"
$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3861"
$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07841"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9605"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.90"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("D12").Value = "1.881.25"
$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.672"
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.881"
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06869"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.24"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009918"
$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.68"
$ws.Range("E19").Value = "  -2.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "28.078.09"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("E22").Value = "  -2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.080"
$ws.Range("E24").Value = "  -3.83%  "

$ws.Range("D25").Value = "2.078.48"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.60"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.700"
$ws.Range("E28").Value = "  -6.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.963"
$ws.Range("E29").Value = "  -3.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.47"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9382"
$ws.Range("E31").Value = "  -3.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09237"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.258"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.319"
$ws.Range("E34").Value = "  -2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.315"
$ws.Range("E35").Value = "  -4.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05834"
$ws.Range("E36").Value = "  -5.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02111"
$ws.Range("E37").Value = "  -4.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.134"
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5579"
$ws.Range("E40").Value = "  -2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.891"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1758"
$ws.Range("E42").Value = "  -2.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07391"
$ws.Range("E43").Value = "  +3.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.63"
$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5262"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.132"
$ws.Range("E46").Value = "  -9.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.096"
$ws.Range("E47").Value = "  -11.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.829"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.01"
$ws.Range("E49").Value = "  -1.82%  "

$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.023"
$ws.Range("E51").Value = "  +0.08%  "
